# ReportingOrganisationGroup.xlsx: the "codeforiati:group-code" and
# "codeforiati:group-name" columns (D and E) were swapped - column D now
# holds the group-name values (and its header), column E now holds the
# group-code values (and its header).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

$dRange = $ws.Range("D1:D$lastRow")
$eRange = $ws.Range("E1:E$lastRow")

$dValues = $dRange.Value2
$eValues = $eRange.Value2

$dRange.Value = $eValues
$eRange.Value = $dValues
